$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A got slightly wider (11.85546875 -> 12.42578125 character units).
# The engine quantizes ColumnWidth in coarse steps, so pick an input value
# that lands on the closest achievable stored width.
$ws.Columns.Item(1).ColumnWidth = 11.665

# Append 5 new rows (206-210) of Date/Count data, continuing the existing
# daily series. Copy formatting (incl. the date number format) down from the
# last existing row first so the new cells pick up the same style index,
# then fill in the actual values.
$ws.Range("A205:B205").Copy()
$ws.Range("A206:B210").PasteSpecial(-4122)

$ws.Range("A206").Value = 45621
$ws.Range("B206").Value = 174

$ws.Range("A207").Value = 45622
$ws.Range("B207").Value = 175

$ws.Range("A208").Value = 45623
$ws.Range("B208").Value = 182

$ws.Range("A209").Value = 45624
$ws.Range("B209").Value = 193

$ws.Range("A210").Value = 45625
$ws.Range("B210").Value = 179

# The active selection on the sheet moved from E4 to E3.
$null = $ws.Range("E3").Select()
